# Actividad 01 - Identificación Hallazgos BD
# Commit: "The excel has been changed, a couple 'Discoveries' have been added"
#
# This script:
#  1. Fixes the titles of Hallazgo 01-03 (adds spaced " - " dash).
#  2. Writes out the full title for Hallazgo 04 and fills in its "Criterio".
#  3. Writes out the full title for Hallazgo 05 and fills in its "Criterio".
#  4. Duplicates the Hallazgo block to create a new "Hallazgo 06" block
#     (rows 31-36) with its "Recomendación" filled in.
#  5. Leaves the "Hallazgos" sheet active, matching the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hallazgos")

# ---------------------------------------------------------------------
# 1. Re-word the existing Hallazgo titles (add spaces around the dash)
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value2 = "Hallazgo 01 - Falta de Restricción de Unicidad en Código de Producto."
$ws.Cells.Item(7, 1).Value2 = "Hallazgo 02 - Productos sin Nombre."
$ws.Cells.Item(13, 1).Value2 = "Hallazgo 03 - Pedidos con Monto Total en Cero."

# ---------------------------------------------------------------------
# 2. Hallazgo 04 - fill in title + Regla de Integridad (Criterio)
# ---------------------------------------------------------------------
$ws.Cells.Item(19, 1).Value2 = "Hallazgo 04 - Llave Foránea no Detectada en Detalle de Pedidos hacia Productos."
$ws.Cells.Item(23, 2).Value2 = "Regla de Integridad Referencial"

# ---------------------------------------------------------------------
# 3. Hallazgo 05 - fill in title + Regla de Integridad (Criterio)
# ---------------------------------------------------------------------
$ws.Cells.Item(25, 1).Value2 = "Hallazgo 05 - Fechas Inconsistentes de Creación"
$ws.Cells.Item(29, 2).Value2 = "Regla de Integridad de Negocio"

# ---------------------------------------------------------------------
# 4. Duplicate the last block (rows 25-30) into a new block for
#    "Hallazgo 06" (rows 31-36), keeping its full formatting/merge,
#    then set its title and Recomendación text.
# ---------------------------------------------------------------------
$ws.Range("A25:B30").Copy($ws.Range("A31"))

$ws.Cells.Item(31, 1).Value2 = "Hallazgo 06"
$ws.Cells.Item(32, 2).Value2 = ""
$ws.Cells.Item(33, 2).Value2 = ""
$ws.Cells.Item(34, 2).Value2 = ""
$ws.Cells.Item(35, 2).Value2 = ""
$ws.Cells.Item(36, 2).Value2 = "Regla de Identidad de Usuario o Dominio"

$ws.Rows.Item(31).RowHeight = 15
$ws.Rows.Item(36).RowHeight = 15

# ---------------------------------------------------------------------
# 5. Leave the view the way it was saved: "Hallazgos" tab active,
#    scrolled down to the new block, B35 selected.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A17").Select()
$ws.Range("B35").Select()
